# Fruta / hortaliza, semanal
# Insert two new weekly records at the top of the data block (rows 236-237),
# pushing the existing rows down by two (236-277 -> 238-279).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 236:277 down to 238:279, leaving two fresh rows at 236:237.
$ws.Rows("236:237").Insert()

# New row 236: Especial quality, week of 2022-02-25 (serial 44617)
$ws.Cells.Item(236, 1).Value  = 7
$ws.Cells.Item(236, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(236, 3).Value  = "Ñuble"
$ws.Cells.Item(236, 4).Value  = 44617
$ws.Cells.Item(236, 5).Value  = 16
$ws.Cells.Item(236, 6).Value  = "Fruta"
$ws.Cells.Item(236, 7).Value  = 100101
$ws.Cells.Item(236, 8).Value  = "Berries"
$ws.Cells.Item(236, 9).Value  = 100112025
$ws.Cells.Item(236, 10).Value = "Frutilla"
$ws.Cells.Item(236, 11).Value = "Sin especificar"
$ws.Cells.Item(236, 12).Value = "Especial"
$ws.Cells.Item(236, 13).Value = 60
$ws.Cells.Item(236, 14).Value = 7000
$ws.Cells.Item(236, 15).Value = 7000
$ws.Cells.Item(236, 16).Value = 7000
$ws.Cells.Item(236, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(236, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(236, 19).Value = 1000
$ws.Cells.Item(236, 20).Value = 7

# New row 237: Primera quality, same week (serial 44617)
$ws.Cells.Item(237, 1).Value  = 7
$ws.Cells.Item(237, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(237, 3).Value  = "Ñuble"
$ws.Cells.Item(237, 4).Value  = 44617
$ws.Cells.Item(237, 5).Value  = 16
$ws.Cells.Item(237, 6).Value  = "Fruta"
$ws.Cells.Item(237, 7).Value  = 100101
$ws.Cells.Item(237, 8).Value  = "Berries"
$ws.Cells.Item(237, 9).Value  = 100112025
$ws.Cells.Item(237, 10).Value = "Frutilla"
$ws.Cells.Item(237, 11).Value = "Sin especificar"
$ws.Cells.Item(237, 12).Value = "Primera"
$ws.Cells.Item(237, 13).Value = 120
$ws.Cells.Item(237, 14).Value = 6000
$ws.Cells.Item(237, 15).Value = 6500
$ws.Cells.Item(237, 16).Value = 6250
$ws.Cells.Item(237, 17).Value = "$/caja 7 kilos"
$ws.Cells.Item(237, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(237, 19).Value = 893
$ws.Cells.Item(237, 20).Value = 7
